# Apply the cryptos-list price/volume update.
# Every text value is written with NumberFormat forced to Text ("@") so
# Excel's automatic type inference can't silently turn number-looking
# strings (e.g. "1.00", "39.60", "59.942.85") into real numbers, then the
# style is reset back to "Normal" so no stray style index is left on the
# cell (matches the source XML, which carries no s=".." on these cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "59.942.85"
Set-TextValue "E2" "  +3.13%  "
Set-TextValue "D3" "2.418.68"
Set-TextValue "E3" "  +2.77%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "552.03"
Set-TextValue "E5" "  +1.86%  "
Set-TextValue "D6" "137.19"
Set-TextValue "E6" "  +2.42%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "E8" "  +3.30%  "
Set-TextValue "E9" "  +1.08%  "
Set-TextValue "D10" "5.76"
Set-TextValue "E10" "  +3.96%  "
Set-TextValue "D11" "0.357"
Set-TextValue "E11" "  +0.09%  "
Set-TextValue "E12" "  -1.84%  "
Set-TextValue "D13" "24.79"
Set-TextValue "E13" "  +4.28%  "
Set-TextValue "D14" "2.850.86"
Set-TextValue "D15" "59.859.98"
Set-TextValue "E15" "  +3.07%  "
Set-TextValue "E16" "  +1.30%  "
Set-TextValue "D17" "2.394.55"
Set-TextValue "E17" "  +2.07%  "
Set-TextValue "D18" "11.31"
Set-TextValue "E18" "  +4.40%  "
Set-TextValue "E19" "  +1.38%  "
Set-TextValue "E20" "  +0.40%  "
Set-TextValue "E21" "  -0.34%  "
Set-TextValue "E22" "  +0.23%  "
Set-TextValue "D23" "65.73"
Set-TextValue "E23" "  +3.77%  "
Set-TextValue "E24" "  +2.85%  "
Set-TextValue "D25" "8.62"
Set-TextValue "E25" "  +4.32%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.12%  "
Set-TextValue "D27" "1.34"
Set-TextValue "E27" "  +1.74%  "
Set-TextValue "D28" "0.0₃0781"
Set-TextValue "E28" "  +6.11%  "
Set-TextValue "E29" "  +0.75%  "
Set-TextValue "D30" "170.21"
Set-TextValue "E30" "  -0.06%  "
Set-TextValue "E31" "  +1.59%  "
Set-TextValue "D32" "18.62"
Set-TextValue "E32" "  +1.25%  "
Set-TextValue "E33" "  +1.70%  "
Set-TextValue "D34" "0.999"
Set-TextValue "E34" "  +0.01%  "
Set-TextValue "E35" "  +4.48%  "
Set-TextValue "E36" "  +0.13%  "
Set-TextValue "E37" "  +0.03%  "
Set-TextValue "E38" "  +0.87%  "
Set-TextValue "D39" "39.60"
Set-TextValue "E40" "  +8.80%  "
Set-TextValue "D41" "312.70"
Set-TextValue "E41" "  +7.95%  "
Set-TextValue "E42" "  +0.65%  "
Set-TextValue "D43" "139.28"
Set-TextValue "E43" "  -1.04%  "
Set-TextValue "E44" "  +1.38%  "
Set-TextValue "D45" "0.0519"
Set-TextValue "E45" "  +1.04%  "
Set-TextValue "B46" "Polygon"
Set-TextValue "C46" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D46" "0.411"
Set-TextValue "E46" "  +7.51%  "
Set-TextValue "D47" "0.577"
Set-TextValue "E47" "  +1.72%  "
Set-TextValue "B48" "InjectiveProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "19.24"
Set-TextValue "E48" "  +1.40%  "
Set-TextValue "E49" "  +1.10%  "
Set-TextValue "D50" "17.63"
Set-TextValue "E50" "  +0.83%  "
Set-TextValue "E51" "  -0.30%  "
